# ETL stage x historico
# The first data row (row 2, sale dated 2022-02-19 / 44611) was removed from
# the staged export; every subsequent row shifts up by one position and the
# trailing row (the old row 27) disappears, shrinking the used range from
# A1:P27 to A1:P26. Deleting the whole sheet row reproduces that shift
# (values, shared string usage and the K-column profit formulas all move
# together), matching the diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete() | Out-Null

# Leave the selection where the editor ended up: the last cell of the
# now-shorter table.
$ws.Range("P26").Select() | Out-Null
